$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "23/10/2025"
$ws.Range("B7").Value = "Al Riyadh"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = "Al Kholood"
$ws.Range("F7").Value = "L"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1.09
$ws.Range("L7").Value = 1.09
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = 14
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 4
